# commit: fix(publipostage): Refactor synthetic array /3
# The "statut" legend used emoji squares (black/green/orange) paired with
# French color-name labels. Swap the black-square entry for a blue-square
# one: ⬛ -> 📘 and its label "noir" -> "bleu". The green/orange pairs keep
# their color-word labels but the square glyphs are also refreshed to the
# new book-emoji set (🟩 -> 📗, 🟧 -> 📙) for visual consistency.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Replace whole-cell matches only, so we don't accidentally touch any
# substring occurrences elsewhere in the sheet.
$used.Replace("⬛", "📘", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$used.Replace("🟩", "📗", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$used.Replace("🟧", "📙", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$used.Replace("noir", "bleu", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
